# Recompute column H ("客単価") on the "ABC分析_客構成" sheet so that it
# reflects the average payment amount (G, "平均支払額") divided by the
# customer count (C, "H.客数（合計）") instead of the previous D.価格 (B) / C.
# Rows where C (customer count) is 0 are left untouched to avoid a
# divide-by-zero (that source row already stores a non-numeric "inf" marker).
# Cells whose recomputed value is unchanged (within floating point noise)
# are left untouched as well, so only the rows that genuinely move get a
# new stored value.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ABC分析_客構成")

$lastRow = $ws.Cells.Item(1, 1).End(-4121).Row
if ($lastRow -lt 2) {
    $lastRow = 74
}

for ($r = 2; $r -le $lastRow; $r++) {
    $c = $ws.Cells.Item($r, 3).Value2
    if ($c -ne 0) {
        $g = $ws.Cells.Item($r, 7).Value2
        $newValue = $g / $c
        $oldValue = $ws.Cells.Item($r, 8).Value2
        $delta = [Math]::Abs($newValue - $oldValue)
        $tolerance = [Math]::Abs($newValue) * 0.0000001
        if ($delta -gt $tolerance) {
            $ws.Cells.Item($r, 8).Value = $newValue
        }
    }
}
